# Update Name of Algo
# Apply updated RandomForest-imputed values to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.835299999999997
$ws.Range("B6").Value = 6.144599999999998
$ws.Range("B7").Value = 5.269199999999999
$ws.Range("C7").Value = -14.16989999999999
$ws.Range("B8").Value = 7.112000000000003
$ws.Range("C11").Value = -12.1937
$ws.Range("C12").Value = -10.8229
$ws.Range("C15").Value = -14.76149999999999
$ws.Range("B16").Value = 7.004099999999998
$ws.Range("B20").Value = 9.163800000000004
$ws.Range("C20").Value = -11.9761
$ws.Range("B21").Value = 9.302000000000001
$ws.Range("C21").Value = -12.08390000000001
$ws.Range("C22").Value = -12.9735
$ws.Range("C23").Value = -11.8281
$ws.Range("B28").Value = 5.619299999999999
$ws.Range("B29").Value = 5.162200000000004
$ws.Range("C29").Value = -10.8449
$ws.Range("B30").Value = 5.0586
$ws.Range("B32").Value = 7.70699999999999
$ws.Range("C34").Value = -11.27490000000001
$ws.Range("B40").Value = 9.355399999999985
$ws.Range("C42").Value = -12.5555
$ws.Range("C43").Value = -13.33449999999999
$ws.Range("C44").Value = -13.84589999999999
$ws.Range("C45").Value = -13.8939
$ws.Range("B46").Value = 6.027399999999997
$ws.Range("C46").Value = -13.634
$ws.Range("C50").Value = -13.95869999999999
$ws.Range("B51").Value = 6.245500000000001
$ws.Range("C51").Value = -12.13
$ws.Range("B52").Value = 5.465600000000002
$ws.Range("B57").Value = 5.672999999999999
$ws.Range("C57").Value = -13.79959999999999
$ws.Range("B59").Value = 5.047299999999998
$ws.Range("B62").Value = 6.071300000000001
$ws.Range("C65").Value = -12.61249999999999
$ws.Range("B66").Value = 5.749600000000003
$ws.Range("C66").Value = -11.5037
$ws.Range("C67").Value = -11.8041
$ws.Range("B73").Value = 8.634299999999998
$ws.Range("B74").Value = 9.144999999999992
$ws.Range("B77").Value = 8.860400000000004
$ws.Range("C79").Value = -11.23750000000001
$ws.Range("C84").Value = -13.26889999999999
$ws.Range("C87").Value = -13.97719999999999
$ws.Range("B92").Value = 5.82649999999999
$ws.Range("C92").Value = -11.4892
$ws.Range("C97").Value = -12.0495
$ws.Range("B100").Value = 5.320699999999996
